$wb = $excel.ActiveWorkbook

# --- Sheet "nulos_por_campo": update a few counts ---
$ws1 = $wb.Worksheets.Item("nulos_por_campo")
$ws1.Range("B6").Value = 227
$ws1.Range("B7").Value = 201
$ws1.Range("B8").Value = 1200

# --- Sheet "quarantine_resumen": add new summary rows ---
$ws2 = $wb.Worksheets.Item("quarantine_resumen")
$ws2.Range("A2").Value = "fecha_invalida"
$ws2.Range("B2").Value = 2906
$ws2.Range("A3").Value = "fecha_invalidasatisf_fuera_rango"
$ws2.Range("B3").Value = 165
$ws2.Range("A4").Value = "satisf_fuera_rango"
$ws2.Range("B4").Value = 91
